$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shift the "MEC-3B-Elem. Máquinas" entry from column C to column F
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "MEC-3B-Elem. Máquinas"

# Row 8: shift the "MEC-3B-Elem. Máquinas" entry from column C to column F
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "MEC-3B-Elem. Máquinas"

# Row 11: shift the "MEC-3A-Metrologia 2" entry from column E to column F
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "[-, -, 'MEC-3A-Metrologia 2', -]"

# Row 12: replace E12 entry with "-" and set F12 to the "MEC-3A-Metrologia 2" entry
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "[-, -, 'MEC-3A-Metrologia 2', -]"

# Row 14: shift the "MEC-3A-Metrologia 2" entry from column E to column F
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "[-, -, 'MEC-3A-Metrologia 2', -]"

# Row 15: shift the "MEC-3A-Metrologia 2" entry from column E to column F
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[-, -, 'MEC-3A-Metrologia 2', -]"
